$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four "Resolving-Mac" sending-cluster rows (old rows 6-9);
# the data now only covers FAPs as the sending cluster (rows 2-5).
$ws.Rows("6:9").Delete()

# Row 2 (FAPs / Rspo1 / Lgr4 / ECs) - refreshed TPM-derived stats
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 1.484826
$ws.Range("N2").Value = 4.454478
$ws.Range("O2").Value = 0.06049021884829667
$ws.Range("P2").Value = 0.06049021884829667
$ws.Range("Q2").Value = 0.173854811746
$ws.Range("R2").Value = 1.564693305714
$ws.Range("S2").Value = 0.06049021884829667
$ws.Range("T2").Value = 0.06049021884829667

# Row 3 (FAPs / Rspo1 / Lgr4 / FAPs) - refreshed TPM-derived stats
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.5859425360316464
$ws.Range("P3").Value = 0.5859425360316464
$ws.Range("S3").Value = 0.5859425360316464
$ws.Range("T3").Value = 0.5859425360316464

# Row 4 (FAPs / Rspo1 / Lgr4 / MuSCs) - refreshed TPM-derived stats
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 8.653369666666666
$ws.Range("N4").Value = 25.960109
$ws.Range("O4").Value = 0.3525289999716321
$ws.Range("P4").Value = 0.3525289999716321
$ws.Range("Q4").Value = 1.013202863074111
$ws.Range("R4").Value = 9.118825767666999
$ws.Range("S4").Value = 0.3525289999716321
$ws.Range("T4").Value = 0.3525289999716321

# Row 5 (FAPs / Rspo1 / Lgr4 / Resolving-Mac) - refreshed TPM-derived stats
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("M5").Value = 0.02548533333333333
$ws.Range("N5").Value = 0.076456
$ws.Range("O5").Value = 0.001038245148424882
$ws.Range("P5").Value = 0.001038245148424882
$ws.Range("Q5").Value = 0.002984018214222222
$ws.Range("R5").Value = 0.026856163928
$ws.Range("S5").Value = 0.001038245148424882
$ws.Range("T5").Value = 0.001038245148424882
